$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.026992
$ws.Range("H2").Value = 0.08097599999999999
$ws.Range("I2").Value = 0.004182906599909731
$ws.Range("J2").Value = 0.00420788870005516
$ws.Range("M2").Value = 0.169654
$ws.Range("N2").Value = 0.508962
$ws.Range("O2").Value = 0.006094264463659866
$ws.Range("P2").Value = 0.006534681579452628
$ws.Range("Q2").Value = 0.004579300768
$ws.Range("R2").Value = 0.041213706912
$ws.Range("S2").Value = 0.00002549173904663819
$ws.Range("T2").Value = 0.00002749721277663732
$ws.Range("G3").Value = 0.026992
$ws.Range("H3").Value = 0.08097599999999999
$ws.Range("I3").Value = 0.004182906599909731
$ws.Range("J3").Value = 0.00420788870005516
$ws.Range("O3").Value = 0.7900017288527916
$ws.Range("P3").Value = 0.8470931604713817
$ws.Range("Q3").Value = 0.59361643152
$ws.Range("R3").Value = 5.342547883679999
$ws.Range("S3").Value = 0.003304503445558439
$ws.Range("T3").Value = 0.00356447373784154
$ws.Range("G4").Value = 0.026992
$ws.Range("H4").Value = 0.08097599999999999
$ws.Range("I4").Value = 0.004182906599909731
$ws.Range("J4").Value = 0.00420788870005516
$ws.Range("M4").Value = 0.04769766666666667
$ws.Range("N4").Value = 0.143093
$ws.Range("O4").Value = 0.001713382501834088
$ws.Range("P4").Value = 0.001837204332049573
$ws.Range("Q4").Value = 0.001287455418666667
$ws.Range("R4").Value = 0.011587098768
$ws.Range("S4").Value = 0.000007166918975091654
$ws.Range("T4").Value = 0.000007730751348523786
$ws.Range("G5").Value = 0.026992
$ws.Range("H5").Value = 0.08097599999999999
$ws.Range("I5").Value = 0.004182906599909731
$ws.Range("J5").Value = 0.00420788870005516
$ws.Range("M5").Value = 5.6286445
$ws.Range("N5").Value = 11.257289
$ws.Range("O5").Value = 0.2021906241817143
$ws.Range("P5").Value = 0.1445349536171162
$ws.Range("Q5").Value = 0.151928372344
$ws.Range("R5").Value = 0.911570234064
$ws.Range("S5").Value = 0.0008457444963295607
$ws.Range("T5").Value = 0.0006081869980884599
$ws.Range("I6").Value = 0.976381346197431
$ws.Range("J6").Value = 0.9822127115383066
$ws.Range("M6").Value = 0.169654
$ws.Range("N6").Value = 0.508962
$ws.Range("O6").Value = 0.006094264463659866
$ws.Range("P6").Value = 0.006534681579452628
$ws.Range("Q6").Value = 1.068908363528667
$ws.Range("R6").Value = 9.620175271758001
$ws.Range("S6").Value = 0.005950326141111384
$ws.Range("T6").Value = 0.00641844731319359
$ws.Range("I7").Value = 0.976381346197431
$ws.Range("J7").Value = 0.9822127115383066
$ws.Range("O7").Value = 0.7900017288527916
$ws.Range("P7").Value = 0.8470931604713817
$ws.Range("S7").Value = 0.7713429515155865
$ws.Range("T7").Value = 0.8320256700721497
$ws.Range("I8").Value = 0.976381346197431
$ws.Range("J8").Value = 0.9822127115383066
$ws.Range("M8").Value = 0.04769766666666667
$ws.Range("N8").Value = 0.143093
$ws.Range("O8").Value = 0.001713382501834088
$ws.Range("P8").Value = 0.001837204332049573
$ws.Range("Q8").Value = 0.3005200868874444
$ws.Range("R8").Value = 2.704680781987
$ws.Range("S8").Value = 0.001672914713691889
$ws.Range("T8").Value = 0.001804525448632334
$ws.Range("I9").Value = 0.976381346197431
$ws.Range("J9").Value = 0.9822127115383066
$ws.Range("M9").Value = 5.6286445
$ws.Range("N9").Value = 11.257289
$ws.Range("O9").Value = 0.2021906241817143
$ws.Range("P9").Value = 0.1445349536171162
$ws.Range("Q9").Value = 35.46338536892517
$ws.Range("R9").Value = 212.780312213551
$ws.Range("S9").Value = 0.1974151538270411
$ws.Range("T9").Value = 0.1419640687043311
$ws.Range("G10").Value = 0.1149325
$ws.Range("H10").Value = 0.229865
$ws.Range("I10").Value = 0.01781090370458377
$ws.Range("J10").Value = 0.01194485200600399
$ws.Range("M10").Value = 0.169654
$ws.Range("N10").Value = 0.508962
$ws.Range("O10").Value = 0.006094264463659866
$ws.Range("P10").Value = 0.006534681579452628
$ws.Range("Q10").Value = 0.019498758355
$ws.Range("R10").Value = 0.11699255013
$ws.Range("S10").Value = 0.0001085443575125127
$ws.Range("T10").Value = 0.00007805580437292207
$ws.Range("G11").Value = 0.1149325
$ws.Range("H11").Value = 0.229865
$ws.Range("I11").Value = 0.01781090370458377
$ws.Range("J11").Value = 0.01194485200600399
$ws.Range("O11").Value = 0.7900017288527916
$ws.Range("P11").Value = 0.8470931604713817
$ws.Range("Q11").Value = 2.527631169075
$ws.Range("R11").Value = 15.16578701445
$ws.Range("S11").Value = 0.01407064471905177
$ws.Range("T11").Value = 0.01011840243712885
$ws.Range("G12").Value = 0.1149325
$ws.Range("H12").Value = 0.229865
$ws.Range("I12").Value = 0.01781090370458377
$ws.Range("J12").Value = 0.01194485200600399
$ws.Range("M12").Value = 0.04769766666666667
$ws.Range("N12").Value = 0.143093
$ws.Range("O12").Value = 0.001713382501834088
$ws.Range("P12").Value = 0.001837204332049573
$ws.Range("Q12").Value = 0.005482012074166666
$ws.Range("R12").Value = 0.032892072445
$ws.Range("S12").Value = 0.00003051689074928577
$ws.Range("T12").Value = 0.00002194513385112157
$ws.Range("G13").Value = 0.1149325
$ws.Range("H13").Value = 0.229865
$ws.Range("I13").Value = 0.01781090370458377
$ws.Range("J13").Value = 0.01194485200600399
$ws.Range("M13").Value = 5.6286445
$ws.Range("N13").Value = 11.257289
$ws.Range("O13").Value = 0.2021906241817143
$ws.Range("P13").Value = 0.1445349536171162
$ws.Range("Q13").Value = 0.64691418399625
$ws.Range("R13").Value = 2.587656735985
$ws.Range("S13").Value = 0.0036011977372702
$ws.Range("T13").Value = 0.001726448630651104
$ws.Range("G14").Value = 0.010485
$ws.Range("H14").Value = 0.031455
$ws.Range("I14").Value = 0.001624843498075486
$ws.Range("J14").Value = 0.001634547755634201
$ws.Range("M14").Value = 0.169654
$ws.Range("N14").Value = 0.508962
$ws.Range("O14").Value = 0.006094264463659866
$ws.Range("P14").Value = 0.006534681579452628
$ws.Range("Q14").Value = 0.00177882219
$ws.Range("R14").Value = 0.01600939971
$ws.Range("S14").Value = 0.000009902225989330225
$ws.Range("T14").Value = 0.00001068124910947845
$ws.Range("G15").Value = 0.010485
$ws.Range("H15").Value = 0.031455
$ws.Range("I15").Value = 0.001624843498075486
$ws.Range("J15").Value = 0.001634547755634201
$ws.Range("O15").Value = 0.7900017288527916
$ws.Range("P15").Value = 0.8470931604713817
$ws.Range("Q15").Value = 0.23058937035
$ws.Range("R15").Value = 2.07530433315
$ws.Range("S15").Value = 0.001283629172594852
$ws.Range("T15").Value = 0.001384614224261579
$ws.Range("G16").Value = 0.010485
$ws.Range("H16").Value = 0.031455
$ws.Range("I16").Value = 0.001624843498075486
$ws.Range("J16").Value = 0.001634547755634201
$ws.Range("M16").Value = 0.04769766666666667
$ws.Range("N16").Value = 0.143093
$ws.Range("O16").Value = 0.001713382501834088
$ws.Range("P16").Value = 0.001837204332049573
$ws.Range("Q16").Value = 0.000500110035
$ws.Range("R16").Value = 0.004500990315
$ws.Range("S16").Value = 0.000002783978417821428
$ws.Range("T16").Value = 0.000003002998217593061
$ws.Range("G17").Value = 0.010485
$ws.Range("H17").Value = 0.031455
$ws.Range("I17").Value = 0.001624843498075486
$ws.Range("J17").Value = 0.001634547755634201
$ws.Range("M17").Value = 5.6286445
$ws.Range("N17").Value = 11.257289
$ws.Range("O17").Value = 0.2021906241817143
$ws.Range("P17").Value = 0.1445349536171162
$ws.Range("Q17").Value = 0.05901633758249999
$ws.Range("R17").Value = 0.354098025495
$ws.Range("S17").Value = 0.0003285281210734827
$ws.Range("T17").Value = 0.0002362492840455506
